# Fix typo "Manualy" -> "Manually" and "selct" -> "select" in the
# ATDD Scenarios worksheet (Table2, "Given-When-Then (Description)" column, G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

$oldGiven = "Manualy created warehouse shipment from released sales order with one line with require shipment location"
$newGiven = "Manually created warehouse shipment from released sales order with one line with require shipment location"

$oldWhen = "Delete warehouse shipment line and selct yes in confirmation"
$newWhen = "Delete warehouse shipment line and select yes in confirmation"

$ws.Range("G63").Value = $newWhen

foreach ($cell in @("G7", "G21", "G41", "G55")) {
    $ws.Range($cell).Value = $newGiven
}

# Update sheet view / window state to match the saved workbook.
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("G7").Select()
